$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 17858914
$ws.Range("I28").Value = 20835192
$ws.Range("J28").Value = 1252.25
$ws.Range("K28").Value = 20835192
$ws.Range("L28").Value = 1252.25
$ws.Range("M28").Value = -20834707
$ws.Range("N28").Value = -2222.25
$ws.Range("H41").Value = 881.7368
$ws.Range("I41").Value = 132.33333
$ws.Range("J41").Value = 1556.2
$ws.Range("K41").Value = 132.33333
$ws.Range("L41").Value = 1556.2
$ws.Range("M41").Value = 307.66667
$ws.Range("N41").Value = -2436.2
$ws.Range("H62").Value = 67310880
$ws.Range("I62").Value = 22730952
$ws.Range("K62").Value = 22730952
$ws.Range("M62").Value = -22730328
$ws.Range("H65").Value = 67310880
$ws.Range("I65").Value = 22730952
$ws.Range("K65").Value = 113654760
$ws.Range("M65").Value = -113651640
$ws.Range("H76").Value = 66670580
$ws.Range("I76").Value = 90913180
$ws.Range("J76").Value = 3412.5
$ws.Range("K76").Value = 90913180
$ws.Range("L76").Value = 3412.5
$ws.Range("M76").Value = -90912865
$ws.Range("N76").Value = -4042.5
$ws.Range("H79").Value = 66670580
$ws.Range("I79").Value = 90913180
$ws.Range("J79").Value = 3412.5
$ws.Range("K79").Value = 90913180
$ws.Range("L79").Value = 3412.5
$ws.Range("M79").Value = -90912088
$ws.Range("N79").Value = -5596.5
$ws.Range("H92").Value = 794.14813
$ws.Range("I92").Value = 873.7826
$ws.Range("J92").Value = 336.25
$ws.Range("K92").Value = 873.7826
$ws.Range("L92").Value = 336.25
$ws.Range("M92").Value = 374.2174
$ws.Range("N92").Value = -2832.25
$ws.Range("H98").Value = 118584870
$ws.Range("I98").Value = 50003748
$ws.Range("J98").Value = 255747120
$ws.Range("K98").Value = 50003748
$ws.Range("L98").Value = 255747120
$ws.Range("M98").Value = -50002250
$ws.Range("N98").Value = -255750116
$ws.Range("H106").Value = 58825950
$ws.Range("I106").Value = 58825950
$ws.Range("K106").Value = 58825950
$ws.Range("M106").Value = -58825319
$ws.Range("H107").Value = 765.36664
$ws.Range("I107").Value = 878.375
$ws.Range("J107").Value = 313.33334
$ws.Range("K107").Value = 878.375
$ws.Range("L107").Value = 313.33334
$ws.Range("M107").Value = 1041.625
$ws.Range("N107").Value = -4153.33334
$ws.Range("H122").Value = 118584870
$ws.Range("I122").Value = 50003748
$ws.Range("J122").Value = 255747120
$ws.Range("K122").Value = 150011244
$ws.Range("L122").Value = 767241360
$ws.Range("M122").Value = -150008794
$ws.Range("N122").Value = -767246260

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10101.15
$ws.Range("I32").Value = 7158.5117
$ws.Range("J32").Value = 25550
$ws.Range("K32").Value = 7158.5117
$ws.Range("L32").Value = 25550
$ws.Range("M32").Value = -6871.5117
$ws.Range("N32").Value = -26124
$ws.Range("H37").Value = 9974.947
$ws.Range("J37").Value = 14152.4
$ws.Range("L37").Value = 14152.4
$ws.Range("N37").Value = -14698.4
$ws.Range("H63").Value = 2068.75
$ws.Range("I63").Value = 2045.4546
$ws.Range("J63").Value = 2120
$ws.Range("K63").Value = 2045.4546
$ws.Range("L63").Value = 2120
$ws.Range("M63").Value = -1359.4546
$ws.Range("N63").Value = -3492
$ws.Range("H66").Value = 2068.75
$ws.Range("I66").Value = 2045.4546
$ws.Range("J66").Value = 2120
$ws.Range("K66").Value = 10227.273
$ws.Range("L66").Value = 10600
$ws.Range("M66").Value = -6795.273000000001
$ws.Range("N66").Value = -17464

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1596.8125
$ws.Range("I105").Value = 1599.9231
$ws.Range("J105").Value = 1583.3334
$ws.Range("K105").Value = 1599.9231
$ws.Range("L105").Value = 1583.3334
$ws.Range("M105").Value = 147.0769
$ws.Range("N105").Value = -5077.3334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H51").Value = 9173.666999999999
$ws.Range("J51").Value = 9173.666999999999
$ws.Range("L51").Value = 9173.666999999999
$ws.Range("N51").Value = -10645.667
$ws.Range("H59").Value = 17825.75
$ws.Range("J59").Value = 17825.75
$ws.Range("L59").Value = 17825.75
$ws.Range("N59").Value = -20115.75
$ws.Range("H60").Value = 8400.666999999999
$ws.Range("J60").Value = 8400.666999999999
$ws.Range("L60").Value = 8400.666999999999
$ws.Range("N60").Value = -9422.666999999999
$ws.Range("H61").Value = 9173.666999999999
$ws.Range("J61").Value = 9173.666999999999
$ws.Range("L61").Value = 9173.666999999999
$ws.Range("N61").Value = -9869.666999999999
$ws.Range("H62").Value = 6759180.5
$ws.Range("I62").Value = 2256.3333
$ws.Range("J62").Value = 35717428
$ws.Range("K62").Value = 2256.3333
$ws.Range("L62").Value = 35717428
$ws.Range("M62").Value = -1632.3333
$ws.Range("N62").Value = -35718676
$ws.Range("H65").Value = 6759180.5
$ws.Range("I65").Value = 2256.3333
$ws.Range("J65").Value = 35717428
$ws.Range("K65").Value = 11281.6665
$ws.Range("L65").Value = 178587140
$ws.Range("M65").Value = -8161.666499999999
$ws.Range("N65").Value = -178593380
$ws.Range("H74").Value = 23942.2
$ws.Range("J74").Value = 23942.2
$ws.Range("L74").Value = 23942.2
$ws.Range("N74").Value = -25690.2
$ws.Range("H77").Value = 23942.2
$ws.Range("J77").Value = 23942.2
$ws.Range("L77").Value = 71826.60000000001
$ws.Range("N77").Value = -80562.60000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 5567637.5
$ws.Range("I131").Value = 45545740
$ws.Range("J131").Value = 1066.2659
$ws.Range("K131").Value = 136637220
$ws.Range("L131").Value = 3198.7977
$ws.Range("M131").Value = -136632180
$ws.Range("N131").Value = -13278.7977
$ws.Range("H137").Value = 5502.9395
$ws.Range("I137").Value = 2040.5625
$ws.Range("J137").Value = 8761.647000000001
$ws.Range("K137").Value = 6121.6875
$ws.Range("L137").Value = 26284.941
$ws.Range("M137").Value = -1021.6875
$ws.Range("N137").Value = -36484.94100000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5138003.5
$ws.Range("I70").Value = 2236393.8
$ws.Range("J70").Value = 11908427
$ws.Range("K70").Value = 2236393.8
$ws.Range("L70").Value = 11908427
$ws.Range("M70").Value = -2236123.8
$ws.Range("N70").Value = -11908967
$ws.Range("H73").Value = 5138003.5
$ws.Range("I73").Value = 2236393.8
$ws.Range("J73").Value = 11908427
$ws.Range("K73").Value = 2236393.8
$ws.Range("L73").Value = 11908427
$ws.Range("M73").Value = -2235457.8
$ws.Range("N73").Value = -11910299
$ws.Range("H80").Value = 8846.429
$ws.Range("I80").Value = 4189.4736
$ws.Range("J80").Value = 18677.777
$ws.Range("K80").Value = 4189.4736
$ws.Range("L80").Value = 18677.777
$ws.Range("M80").Value = -3191.4736
$ws.Range("N80").Value = -20673.777
$ws.Range("H83").Value = 8846.429
$ws.Range("I83").Value = 4189.4736
$ws.Range("J83").Value = 18677.777
$ws.Range("K83").Value = 20947.368
$ws.Range("L83").Value = 93388.88499999999
$ws.Range("M83").Value = -15955.368
$ws.Range("N83").Value = -103372.885

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 4631882
$ws.Range("I136").Value = 5954514.5
$ws.Range("J136").Value = 2667.5
$ws.Range("K136").Value = 17863543.5
$ws.Range("L136").Value = 8002.5
$ws.Range("M136").Value = -17860993.5
$ws.Range("N136").Value = -13102.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1801.0385
$ws.Range("I136").Value = 1453.05
$ws.Range("J136").Value = 2018.5312
$ws.Range("K136").Value = 4359.15
$ws.Range("L136").Value = 6055.5936
$ws.Range("M136").Value = -1809.15
$ws.Range("N136").Value = -11155.5936
